$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text storage for numeric-looking Price/Volume cells so they
# keep matching the original inline-string (text) representation
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","E16","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D38","E38","D39","E39","D40","E40","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "306.96"
$ws.Range("E2").Value = "-6.33%"
$ws.Range("D3").Value = "39.36"
$ws.Range("E3").Value = "-11.39%"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").Value = "-5.40%"
$ws.Range("D5").Value = "0.07735"
$ws.Range("E5").Value = "-7.59%"
$ws.Range("D6").Value = "4.289"
$ws.Range("E6").Value = "-3.11%"
$ws.Range("D7").Value = "1.570"
$ws.Range("E7").Value = "-19.22%"
$ws.Range("D8").Value = "0.9175"
$ws.Range("E8").Value = "-5.54%"
$ws.Range("D9").Value = "0.1022"
$ws.Range("E9").Value = "-9.70%"
$ws.Range("D10").Value = "0.1724"
$ws.Range("E10").Value = "-9.04%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.04477"
$ws.Range("E11").Value = "-2.01%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08906"
$ws.Range("E12").Value = "-7.59%"
$ws.Range("D13").Value = "7.034"
$ws.Range("E13").Value = "-16.02%"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").Value = "-0.12%"
$ws.Range("D15").Value = "0.001280"
$ws.Range("E15").Value = "-0.46%"
$ws.Range("E16").Value = "-1.50%"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("D18").Value = "2.559"
$ws.Range("E18").Value = "1.07%"
$ws.Range("E19").Value = "0.31%"
$ws.Range("D20").Value = "0.1365"
$ws.Range("E20").Value = "1.01%"
$ws.Range("D21").Value = "0.2777"
$ws.Range("E21").Value = "1.89%"
$ws.Range("D22").Value = "0.04132"
$ws.Range("E22").Value = "-0.99%"
$ws.Range("D23").Value = "0.001198"
$ws.Range("E23").Value = "-3.36%"
$ws.Range("D24").Value = "0.004081"
$ws.Range("E24").Value = "-7.84%"
$ws.Range("D25").Value = "0.0001225"
$ws.Range("E25").Value = "-5.93%"
$ws.Range("D26").Value = "0.0002991"
$ws.Range("E26").Value = "0.38%"
$ws.Range("D38").Value = "0.02356"
$ws.Range("E38").Value = "-12.97%"
$ws.Range("D39").Value = "0.05136"
$ws.Range("E39").Value = "-8.60%"
$ws.Range("D40").Value = "0.007956"
$ws.Range("E40").Value = "2.23%"
$ws.Range("E41").Value = "-5.78%"
$ws.Range("D42").Value = "0.007389"
$ws.Range("E42").Value = "0.35%"
$ws.Range("D43").Value = "0.001995"
$ws.Range("E43").Value = "-5.99%"
$ws.Range("D44").Value = "0.008040"
$ws.Range("E44").Value = "-7.60%"
$ws.Range("D45").Value = "0.3322"
$ws.Range("E45").Value = "-5.37%"
$ws.Range("D46").Value = "0.00006692"
$ws.Range("E46").Value = "-2.79%"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").Value = "0.38%"
$ws.Range("D48").Value = "0.003396"
$ws.Range("E48").Value = "-2.92%"
$ws.Range("D49").Value = "0.004117"
$ws.Range("E49").Value = "16.59%"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").Value = "0.38%"
$ws.Range("D51").Value = "0.0002008"
$ws.Range("E51").Value = "0.38%"
